$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("1").Delete()
$ws.Range("A1").ClearFormats()
$newValue = @'
questions = [
    {
        "title": "You are developing a mobile-first application using Bootstrap. The web pages should be displayed using 10 columns on all devices, with a one-column margin on the left side.  Which of the following code snippets should you use to implement this?",
        "ques_type": 2,
        "options": [
            "&ltdiv class=\"row\"&gt&ltdiv class=\"col-1\"&gt&lt/div&gt&ltdiv class=\"col-10\"&gt&lt/div&gt&lt/div&gt",
            "&ltdiv class=\"row\"&gt&ltdiv class=\"col-sm-1\"&gt&lt/div&gt&ltdiv class=\"col-10\"&gt&lt/div&gt&lt/div&gt",
            "&ltdiv class=\"row\"&gt&ltdiv class=\"col-md-1\"&gt&lt/div&gt&ltdiv class=\"col-10\"&gt&lt/div&gt&lt/div&gt",
            "&ltdiv class=\"row\"&gt&ltdiv class=\"col-lg-1\"&gt&lt/div&gt&ltdiv class=\"col-10\"&gt&lt/div&gt&lt/div&gt"
        ],
        "score": "&ltdiv class=\"row\"&gt&ltdiv class=\"col-1\"&gt&lt/div&gt&ltdiv class=\"col-10\"&gt&lt/div&gt&lt/div&gt"
    },
    {
        "title": "You are creating a form for data entry using Bootstrap. You want a form input that displays at the tallest height on all devices, regardless of screen size.  Which form input class should you select?",
        "ques_type": 2,
        "options": [
            "form-control form-control-md",
            "form-control",
            "form-control form-control-lg",
            "form-control form-control-sm"
        ],
        "score": "form-control form-control-lg"
    },
    {
        "title": "You have developed a web page that requires a \u201csubmit\u201d button for submitting a form. You want the button to appear with the default blue color.  Which button class should you select?",
        "ques_type": 2,
        "options": [
            "btn btn-light",
            "btn btn-info",
            "btn btn-success",
            "btn btn-primary"
        ],
        "score": "btn btn-primary"
    },
    {
        "title": "You want to float text on a web page that is displayed using Bootstrap. You want the text to float only on large or extra large devices.  True or false:  You should choose the class float-md-start to achieve this.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "False"
    }
]
'@
$ws.Range("A1").Value = $newValue
$ws.Range("A1").EntireRow.AutoFit()
Write-Host "done"
